# DPLKINV107-001 - Profile Fixed Income - Fixed Income View Detil
# Update the expected "Kode Fixed Income" test data from OBL00107 -> OBL00108,
# reflecting this both in the KODE_FIXED_INCOME cell (M2) and inside the
# PREPARATION instructions cell (F2). Also move the active selection to G2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# KODE_FIXED_INCOME column (M2): OBL00107 -> OBL00108
$ws.Range("M2").Value = "OBL00108"

# PREPARATION column (F2): update the "Kode Fixed Income" line inside the
# multi-line preparation text to reference the new code.
$ws.Range("F2").Value = "Username : 31246;`nPassword : bni1234;`nRole : 20/21 - Analis Investasi/Asisten Investasi;`nKode Fixed Income : OBL00108"

# Move the active selection from N2 to G2.
$ws.Range("G2").Select() | Out-Null
